$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Add the three new header cells (L1:N1) with the new column names.
#    These reuse the same bold header style ("s=3") as the existing headers
#    simply by copying the format from K1 (the last existing header cell).
$ws.Range("K1").Copy() | Out-Null
$ws.Range("L1:N1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# 2) taxa_sucesso (E) and particip (F) switch from fraction (0-1) to
#    percentage-point scale (0-100), keeping their existing "0.00%" style.
$ws.Range("E2").Value = 70.0374531835206
$ws.Range("F2").Value = 62.56684491978609

$ws.Range("E3").Value = 29.9625468164794
$ws.Range("F3").Value = 61.25000000000001

$ws.Range("E4").Value = 61.30790190735694
$ws.Range("F4").Value = 93

$ws.Range("E5").Value = 38.69209809264305
$ws.Range("F5").Value = 96.12676056338029

$ws.Range("E6").Value = 90.93567251461988
$ws.Range("F6").Value = 20.90032154340836

$ws.Range("E7").Value = 9.064327485380117
$ws.Range("F7").Value = 35.48387096774194

# 3) Fill the new L/M/N columns (apoio_medio, contribuicoes, media_contribuicoes)
#    for every data row. These have no special number format (general format).
$ws.Range("L2").Value = 90.81853194977892
$ws.Range("M2").Value = 193026
$ws.Range("N2").Value = 329.9589743589743

$ws.Range("L3").Value = 92.63036679831843
$ws.Range("M3").Value = 70527
$ws.Range("N3").Value = 287.865306122449

$ws.Range("L4").Value = 84.63408307975531
$ws.Range("M4").Value = 126119
$ws.Range("N4").Value = 150.6798088410992

$ws.Range("L5").Value = 99.16759340131101
$ws.Range("M5").Value = 77527
$ws.Range("N5").Value = 141.9908424908425

$ws.Range("L6").Value = 18.15137523021585
$ws.Range("M6").Value = 1885
$ws.Range("N6").Value = 14.5

$ws.Range("L7").Value = 27.77589921308953
$ws.Range("M7").Value = 323
$ws.Range("N7").Value = 14.68181818181818
